$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.458829987332251
$ws.Range("C2").Value = 0.3733463158623067
$ws.Range("D2").Value = 0.07915418184603595
$ws.Range("E2").Value = 0.4194773432817129
$ws.Range("G2").Value = 1.067603345455211
$ws.Range("H2").Value = 0.9205681559414813
$ws.Range("N2").Value = 1.004369710222033
$ws.Range("B3").Value = 1.296233088693612
$ws.Range("C3").Value = 0.3248722604003547
$ws.Range("D3").Value = 0.07167795652539155
$ws.Range("E3").Value = 0.3655159147182587
$ws.Range("G3").Value = 1.007881880447229
$ws.Range("H3").Value = 0.8998350751053579
$ws.Range("N3").Value = 1.020781354355904
$ws.Range("B4").Value = 1.196917738862737
$ws.Range("C4").Value = 0.2951292708929145
$ws.Range("D4").Value = 0.06713178741078707
$ws.Range("E4").Value = 0.3325322338200181
$ws.Range("G4").Value = 0.9721569116691455
$ws.Range("H4").Value = 0.8878196540285046
$ws.Range("N4").Value = 1.031384026589798
$ws.Range("B5").Value = 1.156572705499684
$ws.Range("C5").Value = 0.2830125298288522
$ws.Range("D5").Value = 0.06529007773322348
$ws.Range("E5").Value = 0.3191252144372783
$ws.Range("G5").Value = 0.9578319292870106
$ws.Range("H5").Value = 0.8831012017985813
$ws.Range("N5").Value = 1.035836247133087
$ws.Range("B6").Value = 1.149880987233303
$ws.Range("C6").Value = 0.2810007434346744
$ws.Range("D6").Value = 0.06498491412483531
$ws.Range("E6").Value = 0.3169009520463533
$ws.Range("G6").Value = 0.9554672414624292
$ws.Range("H6").Value = 0.8823284005372898
$ws.Range("N6").Value = 1.03658346010722
$ws.Range("B7").Value = 1.196373123454578
$ws.Range("C7").Value = 0.294965846769685
$ws.Range("D7").Value = 0.06710690568560551
$ws.Range("E7").Value = 0.3323512883776232
$ws.Range("G7").Value = 0.9719627811195153
$ws.Range("H7").Value = 0.8877553012672195
$ws.Range("N7").Value = 1.031443539064099
$ws.Range("B8").Value = 1.402656254550834
$ws.Range("C8").Value = 0.3566274320950242
$ws.Range("D8").Value = 0.07656704514189983
$ws.Range("E8").Value = 0.4008385141661535
$ws.Range("G8").Value = 1.046812838934642
$ws.Range("H8").Value = 0.9132699353826865
$ws.Range("N8").Value = 1.00991889074756
$ws.Range("B9").Value = 1.811501836430011
$ws.Range("C9").Value = 0.4777730332828583
$ws.Range("D9").Value = 0.0954814382823912
$ws.Range("E9").Value = 0.5364807574718782
$ws.Range("G9").Value = 1.201285346262864
$ws.Range("H9").Value = 0.9690614951722694
$ws.Range("N9").Value = 0.9719118615105415
$ws.Range("B10").Value = 2.114828834429886
$ws.Range("C10").Value = 0.5670189753116688
$ws.Range("D10").Value = 0.1096172124639594
$ws.Range("E10").Value = 0.6371832104419468
$ws.Range("G10").Value = 1.319755956591422
$ws.Range("H10").Value = 1.013687192842724
$ws.Range("N10").Value = 0.9465913107503141
$ws.Range("B11").Value = 2.253530972620638
$ws.Range("C11").Value = 0.6076943741458649
$ws.Range("D11").Value = 0.116103968987801
$ws.Range("E11").Value = 0.6832737295401756
$ws.Range("G11").Value = 1.374797109839733
$ws.Range("H11").Value = 1.034806053287554
$ws.Range("N11").Value = 0.9356466480480421
$ws.Range("B12").Value = 2.306162239173432
$ws.Range("C12").Value = 0.6231098880763852
$ws.Range("D12").Value = 0.1185687403045819
$ws.Range("E12").Value = 0.7007714289187987
$ws.Range("G12").Value = 1.395809994045464
$ws.Range("H12").Value = 1.042923130010621
$ws.Range("N12").Value = 0.9315855012788816
$ws.Range("B13").Value = 2.294822285517796
$ws.Range("C13").Value = 0.6197892947658943
$ws.Range("D13").Value = 0.1180375306210095
$ws.Range("E13").Value = 0.6970009567829152
$ws.Range("G13").Value = 1.391276854422813
$ws.Range("H13").Value = 1.04116960868663
$ws.Range("N13").Value = 0.932456420741083
$ws.Range("B14").Value = 2.257858791527838
$ws.Range("C14").Value = 0.6089623557270443
$ws.Range("D14").Value = 0.1163065776970456
$ws.Range("E14").Value = 0.6847123669843995
$ws.Range("G14").Value = 1.376522417318711
$ws.Range("H14").Value = 1.035471436266391
$ws.Range("N14").Value = 0.9353108594843533
$ws.Range("B15").Value = 2.235231759928581
$ws.Range("C15").Value = 0.6023322336354227
$ws.Range("D15").Value = 0.1152474183958816
$ws.Range("E15").Value = 0.6771911265255
$ws.Range("G15").Value = 1.367507182892609
$ws.Range("H15").Value = 1.031996809685438
$ws.Range("N15").Value = 0.9370701670083363
$ws.Range("B16").Value = 2.105779149482032
$ws.Range("C16").Value = 0.5643624232675961
$ws.Range("D16").Value = 0.109194444677712
$ws.Range("E16").Value = 0.6341770609332258
$ws.Range("G16").Value = 1.316182408698126
$ws.Range("H16").Value = 1.012323675795358
$ws.Range("N16").Value = 0.9473181884739859
$ws.Range("B17").Value = 2.026551120454428
$ws.Range("C17").Value = 0.5410899833796634
$ws.Range("D17").Value = 0.1054957630240523
$ws.Range("E17").Value = 0.6078639007685922
$ws.Range("G17").Value = 1.284993864244598
$ws.Range("H17").Value = 1.000465966305825
$ws.Range("N17").Value = 0.9537525769570507
$ws.Range("B18").Value = 1.981048471485281
$ws.Range("C18").Value = 0.5277114412504602
$ws.Range("D18").Value = 0.1033736634580293
$ws.Range("E18").Value = 0.5927554284369734
$ws.Range("G18").Value = 1.267162823287237
$ws.Range("H18").Value = 0.993722627518224
$ws.Range("N18").Value = 0.957507422364241
$ws.Range("B19").Value = 1.96565346136714
$ws.Range("C19").Value = 0.5231828779840271
$ws.Range("D19").Value = 0.1026560562028891
$ws.Range("E19").Value = 0.5876443354253524
$ws.Range("G19").Value = 1.261143904258574
$ws.Range("H19").Value = 0.9914525924265831
$ws.Range("N19").Value = 0.9587879897295686
$ws.Range("B20").Value = 2.034978089300182
$ws.Range("C20").Value = 0.5435666232165204
$ws.Range("D20").Value = 0.1058889450977745
$ws.Range("E20").Value = 0.6106622460676618
$ws.Range("G20").Value = 1.288302747689926
$ws.Range("H20").Value = 1.001720265305522
$ws.Range("N20").Value = 0.9530620350432031
$ws.Range("B21").Value = 2.268712903041489
$ws.Range("C21").Value = 0.6121421350086962
$ws.Range("D21").Value = 0.1168147709222609
$ws.Range("E21").Value = 0.6883205921627393
$ws.Range("G21").Value = 1.380851500299002
$ws.Range("H21").Value = 1.037141858444301
$ws.Range("N21").Value = 0.934470172421122
$ws.Range("B22").Value = 2.422102565488785
$ws.Range("C22").Value = 0.6570345220313811
$ws.Range("D22").Value = 0.1240043929593782
$ws.Range("E22").Value = 0.7393344241266391
$ws.Range("G22").Value = 1.442330638515585
$ws.Range("H22").Value = 1.060991356462978
$ws.Range("N22").Value = 0.9228056045774551
$ws.Range("B23").Value = 2.340176362460852
$ws.Range("C23").Value = 0.6330672859429001
$ws.Range("D23").Value = 0.1201625827396384
$ws.Range("E23").Value = 0.7120823204642477
$ws.Range("G23").Value = 1.409425557016334
$ws.Range("H23").Value = 1.048197714932684
$ws.Range("N23").Value = 0.9289864232216161
$ws.Range("B24").Value = 2.031168111885165
$ws.Range("C24").Value = 0.5424469310572704
$ws.Range("D24").Value = 0.1057111740420709
$ws.Range("E24").Value = 0.6093970542924012
$ws.Range("G24").Value = 1.286806491406196
$ws.Range("H24").Value = 1.001152966967823
$ws.Range("N24").Value = 0.9533740558323771
$ws.Range("B25").Value = 1.700403441628339
$ws.Range("C25").Value = 0.444965208069334
$ws.Range("D25").Value = 0.09032387880350257
$ws.Range("E25").Value = 0.499619096640302
$ws.Range("G25").Value = 1.158643226371225
$ws.Range("H25").Value = 0.9533391626897583
$ws.Range("N25").Value = 0.9817396445230386
